# The paragraph "<id>p028v_1</id>" is currently split across three runs:
#   <id>   (Courier New, color 7f6000)
#   p028v_1 (plain, color 000000)
#   </id>  (Courier New, color 7f6000)
# Collapse them into a single run "<id>p028v_1</id>" carrying the
# formatting of the first run, by doing a Find & Replace over the whole
# run (Word merges matched runs into one run using the first run's
# character formatting when the replacement text is applied).

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "<id>p028v_1</id>",  # FindText
    $true,               # MatchCase
    $false,              # MatchWholeWord
    $false,              # MatchWildcards
    $false,              # MatchSoundsLike
    $false,              # MatchAllWordForms
    $true,                # Forward
    1,                   # Wrap (wdFindContinue)
    $false,              # Format
    "<id>p028v_1</id>",  # ReplaceWith
    2                    # Replace (wdReplaceAll)
) | Out-Null
